$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Niveles" row as done by setting its Status cell (C2) to "ok",
# keeping the existing gray-fill style on that cell.
$ws.Range("C2").Value = "ok"

# Update the active selection to reflect where the edit was made.
$ws.Range("C3").Select()
